# Update LoginData.xlsx "Login" sheet test credentials (commit: "commiting code 4th june")

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Login")

# Row 2 previously held a stale Tetherfi QA login (URL / tetherfiqa / T3th3rf!@12345).
# Replace it with the same OCMUI URL already used on row 3, an Administrator
# login type, and a freshly generated password.
$ws.Range("A2").Value = "https://lab.singtel.tetherfi.cloud:45443/OCMUI"
$ws.Range("B2").Value = "Administrator"
$ws.Range("C2").Value = "pSHS-Iq;DXfKp;dAw;Lfufub&CEL*-tD"

# Move the active selection from C12 to B12.
$ws.Activate()
$ws.Range("B12").Select()
